$d = $word.ActiveDocument

# The document ends with a paragraph reading "Graphic is not optimized for
# low resolution" immediately followed by an otherwise-empty paragraph whose
# only content is the "_GoBack" bookmark. Merge the two paragraphs into one
# (same effect as placing the cursor at the end of the text and pressing
# Delete) so the bookmark ends up inside the now-last paragraph, and the
# now-redundant empty paragraph (with its own indentation/font pPr) goes
# away.
$bm = $d.Bookmarks("_GoBack")
$bmStart = $bm.Range.Start

# Range covering just the paragraph mark immediately preceding the bookmark.
$mark = $d.Range($bmStart - 1, $bmStart)
$mark.Delete()
